# Scheduled-runner data refresh: updates the cached market-board price /
# profit figures (columns H-N) on a handful of rows across the ALC, ARM,
# CRP, CUL, GSM, LTW and WVR sheets. Cells that should no longer carry a
# value are cleared outright rather than being set to 0/blank text, so
# that they round-trip the same way the source workbook does (i.e. the
# cell disappears from the row instead of storing an empty value).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 872679.0600000001
$ws.Range("J17").Value = 872679.0600000001
$ws.Range("L17").Value = 2618037.18
$ws.Range("N17").Value = -2618373.18

$ws.Range("H64").Value = 27780800
$ws.Range("I64").Value = 47621256
$ws.Range("J64").Value = 4160
$ws.Range("K64").Value = 47621256
$ws.Range("L64").Value = 4160
$ws.Range("M64").Value = -47621008
$ws.Range("N64").Value = -4656

$ws.Range("H67").Value = 27780800
$ws.Range("I67").Value = 47621256
$ws.Range("J67").Value = 4160
$ws.Range("K67").Value = 47621256
$ws.Range("L67").Value = 4160
$ws.Range("M67").Value = -47620398
$ws.Range("N67").Value = -5876

$ws.Range("H88").Value = 4920.304
$ws.Range("I88").Value = 549.5
$ws.Range("J88").Value = 7251.4
$ws.Range("K88").Value = 549.5
$ws.Range("L88").Value = 7251.4
$ws.Range("M88").Value = -143.5
$ws.Range("N88").Value = -8063.4

$ws.Range("H91").Value = 4920.304
$ws.Range("I91").Value = 549.5
$ws.Range("J91").Value = 7251.4
$ws.Range("K91").Value = 549.5
$ws.Range("L91").Value = 7251.4
$ws.Range("M91").Value = 854.5
$ws.Range("N91").Value = -10059.4

$ws.Range("H138").Value = 2358.8262
$ws.Range("I138").Value = 1354.48
$ws.Range("J138").Value = 3554.476
$ws.Range("K138").Value = 4063.44
$ws.Range("L138").Value = 10663.428
$ws.Range("M138").Value = 1076.56
$ws.Range("N138").Value = -20943.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("N17").Value = 0

$ws.Range("H104").Value = 26000
$ws.Range("J104").Value = 26000
$ws.Range("L104").Value = 26000
$ws.Range("N104").Value = -32988

$ws.Range("H122").Value = 1589
$ws.Range("I122").Value = 1299.7
$ws.Range("K122").Value = 3899.1
$ws.Range("M122").Value = -1449.1

$ws.Range("H124").Value = 22333.334
$ws.Range("J124").Value = 22333.334
$ws.Range("L124").Value = 22333.334
$ws.Range("N124").Value = -32153.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15141.047
$ws.Range("I31").Value = 1161.4694
$ws.Range("J31").Value = 33654.54
$ws.Range("K31").Value = 1161.4694
$ws.Range("L31").Value = 33654.54
$ws.Range("M31").Value = -866.4694
$ws.Range("N31").Value = -34244.54

$ws.Range("H34").Value = 15141.047
$ws.Range("I34").Value = 1161.4694
$ws.Range("J34").Value = 33654.54
$ws.Range("K34").Value = 1161.4694
$ws.Range("L34").Value = 33654.54
$ws.Range("M34").Value = -959.4694
$ws.Range("N34").Value = -34058.54

$ws.Range("H99").Value = 2107.9473
$ws.Range("I99").Value = 1977.4
$ws.Range("J99").Value = 2597.5
$ws.Range("K99").Value = 1977.4
$ws.Range("L99").Value = 2597.5
$ws.Range("M99").Value = -479.4000000000001
$ws.Range("N99").Value = -5593.5

$ws.Range("H108").Value = 48000
$ws.Range("J108").Value = 48000
$ws.Range("L108").Value = 48000
$ws.Range("N108").Value = -55680

$ws.Range("H122").Value = 1395.25
$ws.Range("I122").Value = 1264.8572
$ws.Range("J122").Value = 1496.6666
$ws.Range("K122").Value = 3794.5716
$ws.Range("L122").Value = 4489.9998
$ws.Range("M122").Value = -1344.5716
$ws.Range("N122").Value = -9389.9998

$ws.Range("H126").Value = 2107.9473
$ws.Range("I126").Value = 1977.4
$ws.Range("J126").Value = 2597.5
$ws.Range("K126").Value = 5932.200000000001
$ws.Range("L126").Value = 7792.5
$ws.Range("M126").Value = -3462.200000000001
$ws.Range("N126").Value = -12732.5

$ws.Range("H134").Value = 4487.846
$ws.Range("I134").Value = 6242.857
$ws.Range("J134").Value = 2440.3333
$ws.Range("K134").Value = 18728.571
$ws.Range("L134").Value = 7320.999899999999
$ws.Range("M134").Value = -16193.571
$ws.Range("N134").Value = -12390.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 400
$ws.Range("I75").Value = 400
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 1200
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0
$ws.Range("M75").Value = -202

$ws.Range("H78").Value = 400
$ws.Range("I78").Value = 400
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 3600
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0
$ws.Range("M78").Value = 1392

$ws.Range("H131").Value = 832.53845
$ws.Range("J131").Value = 957.51166
$ws.Range("L131").Value = 2872.53498
$ws.Range("N131").Value = -12952.53498

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2885.8
$ws.Range("I80").Value = 2170
$ws.Range("J80").Value = 3959.5
$ws.Range("K80").Value = 2170
$ws.Range("L80").Value = 3959.5
$ws.Range("M80").Value = -1172
$ws.Range("N80").Value = -5955.5

$ws.Range("H83").Value = 2885.8
$ws.Range("I83").Value = 2170
$ws.Range("J83").Value = 3959.5
$ws.Range("K83").Value = 10850
$ws.Range("L83").Value = 19797.5
$ws.Range("M83").Value = -5858
$ws.Range("N83").Value = -29781.5

$ws.Range("H102").Value = 2855.9048
$ws.Range("I102").Value = 2980.1667
$ws.Range("J102").Value = 2690.2222
$ws.Range("K102").Value = 2980.1667
$ws.Range("L102").Value = 2690.2222
$ws.Range("M102").Value = -1358.1667
$ws.Range("N102").Value = -5934.2222

$ws.Range("H122").Value = 2328.3809
$ws.Range("I122").Value = 1530.5294
$ws.Range("J122").Value = 5719.25
$ws.Range("K122").Value = 4591.5882
$ws.Range("L122").Value = 17157.75
$ws.Range("M122").Value = -2141.5882
$ws.Range("N122").Value = -22057.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2931.8
$ws.Range("I46").Value = 2914.75
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 2914.75
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -2726.75
$ws.Range("N46").Value = -3376

$ws.Range("H93").Value = 1292.3529
$ws.Range("I93").Value = 797.3333
$ws.Range("J93").Value = 1849.25
$ws.Range("K93").Value = 797.3333
$ws.Range("L93").Value = 1849.25
$ws.Range("M93").Value = 450.6667
$ws.Range("N93").Value = -4345.25

$ws.Range("H100").Value = 1548.125
$ws.Range("I100").Value = 1432.3043
$ws.Range("K100").Value = 1432.3043
$ws.Range("M100").Value = -891.3043

$ws.Range("H127").Value = 37540
$ws.Range("J127").Value = 37540
$ws.Range("L127").Value = 37540
$ws.Range("N127").Value = -47460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = 0

$ws.Range("H81").Value = 1800
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 1800
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H86").Value = 28992.5
$ws.Range("J86").Value = 28992.5
$ws.Range("L86").Value = 28992.5
$ws.Range("N86").Value = -31238.5

$ws.Range("H89").Value = 28992.5
$ws.Range("J89").Value = 28992.5
$ws.Range("L89").Value = 144962.5
$ws.Range("N89").Value = -156194.5

$ws.Range("H96").Value = 1286.2858
$ws.Range("I96").Value = 800
$ws.Range("J96").Value = 1480.8
$ws.Range("K96").Value = 800
$ws.Range("L96").Value = 1480.8
$ws.Range("M96").Value = 573
$ws.Range("N96").Value = -4226.8

$ws.Range("H126").Value = 1089.0869
$ws.Range("I126").Value = 972.5333000000001
$ws.Range("J126").Value = 1307.625
$ws.Range("K126").Value = 2917.5999
$ws.Range("L126").Value = 3922.875
$ws.Range("M126").Value = -447.5999000000002
$ws.Range("N126").Value = -8862.875
